# "code carte arriere avec interruptions et maj excel adress CAN"
#
# Adds the missing "Traction control" CAN message row (B16:D16) to the
# "Carte arriere" (rear board) block of the Bus CAN ID list, and updates
# the active cell selection left by the editor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 was previously an empty placeholder row (only A16/D16 styled).
# Fill in the new Traction control CAN id entry: hex id, its decimal
# HEX2DEC conversion, and the message name.
$ws.Range("B16").Value = 1004
$ws.Range("C16").Value = 4100
$ws.Range("D16").Value = "Traction control"

# Leave the selection where the author last left it when saving.
[void]$ws.Range("F16").Select()
